# Auto-save inventory data - 2025-07-04 06:59:50
# Two new inventory rows were inserted in the sheet (pushing the previous
# "B3" rows down), and a Model value was added to one of the moved rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Move the existing "B3" rows (22-23) down to rows 24-25 ---------
$ws.Range("A24").Value = "B3"
$ws.Range("B24").Value = "Stuff"
$ws.Range("C24").Value = 1

$ws.Range("A25").Value = "B3"
$ws.Range("B25").Value = "Test"
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = "123123"
$ws.Range("E25").Value = "123"

# --- 2. Write the new rows (22-23) for location "D2" --------------------
$ws.Range("A22").Value = "D2"
$ws.Range("B22").Value = "New Item"
$ws.Range("C22").Value = 1
$ws.Range("G22").ClearContents()

$ws.Range("A23").Value = "D2"
$ws.Range("B23").Value = "New Item"
$ws.Range("C23").Value = 1
$ws.Range("E23").Value = "123"

# --- 3. Re-point the image hyperlink that used to sit on row 22 to the
#        row it now lives on (row 24). The engine only exposes a
#        sheet-wide Hyperlinks.Delete(), so every hyperlink is rebuilt in
#        its final position/order (same targets, same relative order,
#        with the last one now anchored on G24 instead of G22).
$targets = @(
    "https://www.saturdayeveningpost.com/wp-content/uploads/satevepost/2019-12-19-random-stuff-860x573.jpg",
    "https://media-ecn.s3.amazonaws.com/embedded_image/2016/02/fda.jpg",
    "https://products.integralife.com/ccstore/v1/images/?source=/file/products/Omni-Tract%20Flexible%20Wishbone%20Urologic%20Surgery%20Retractor%20System%20OS%201%20Image.png",
    "https://www.bnsmed.com/data/watermark/20200924/5f6c31aea1382.jpg",
    "https://products.integralife.com/ccstore/v1/images/?source=/file/products/Codman%20Electrosurgical%20Generator%20OS%201%20Image.jpg",
    "https://www.elliquence.com/wp-content/uploads/2016/01/Surgi-Max-Plus-Device.jpg",
    "https://www.aamedicalstore.com/SSP%20Applications/AA%20Medical%20SCA/AA%20Medical/img/Product%20Images/Integra-Duo-LED-Headlight-Set_01.JPG",
    "https://products.integralife.com/ccstore/v1/images/?source=/file/v6400991064904479991/products/MLX-300-Xenon-Lightsources.jpg",
    "https://www.bnsmed.com/data/watermark/20200924/5f6c30bda627b.jpg",
    "https://www.bnsmed.com/data/watermark/20200924/5f6c31aea1382.jpg",
    "https://products.integralife.com/ccstore/v1/images/?source=/file/v3841902670343812321/products/ETK_01.png",
    "https://products.integralife.com/ccstore/v1/images/?source=/file/v3841902670343812321/products/ETK_01.png",
    "https://products.integralife.com/ccstore/v1/images/?source=/file/v7357354864197611707/collections/licox.jpg",
    "https://products.integralife.com/ccstore/v1/images/?source=/file/v5137398853523069574/products/823190.jpg",
    "https://products.integralife.com/ccstore/v1/images/?source=/file/v5137398853523069574/products/823190.jpg",
    "https://products.integralife.com/ccstore/v1/images/?source=/file/v7357354864197611707/collections/licox.jpg",
    "https://products.integralife.com/ccstore/v1/images/?source=/file/v6400991064904479991/products/MLX-300-Xenon-Lightsources.jpg",
    "https://products.integralife.com/ccstore/v1/images/?source=/file/v6400991064904479991/products/MLX-300-Xenon-Lightsources.jpg",
    "https://www.saturdayeveningpost.com/wp-content/uploads/satevepost/2019-12-19-random-stuff-860x573.jpg"
)
$cells = @("G4","G5","G6","G7","G8","G9","G10","G11","G12","G13","G14","G15","G16","G17","G18","G19","G20","G21","G24")

$ws.Range("G4").Hyperlinks.Delete()

for ($i = 0; $i -lt $cells.Count; $i++) {
    $ws.Hyperlinks.Add($ws.Range($cells[$i]), $targets[$i])
}

# --- 4. Dimension updates automatically to A1:G25 based on used range ---
